# Add new columns I (I0) and J (IF) to the sheet, matching the style
# of the existing header cells and filling in the numeric data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style/format from the existing H1 header cell onto
# the two new header cells so they reuse the same (bold/bordered) style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I and J, rows 2-13
$data = @(
    @(5, 7),
    @(8, 9),
    @(3, 6),
    @(9, 9),
    @(7, 8),
    @(7, 9),
    @(6, 7),
    @(8, 8),
    @(2, 8),
    @(4, 5),
    @(5, 6),
    @(5, 5)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
